# Update cryptocurrency price/volume data (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.502.82"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.011.81"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'260.91"
$ws.Range("E5").Value = "  +5.30%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'55.54"
$ws.Range("E8").Value = "  -8.29%  "
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'0.0772"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "2.311.42"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "'14.30"
$ws.Range("E13").Value = "  -5.50%  "
$ws.Range("D14").Value = "'0.801"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").Value = "'20.63"
$ws.Range("E15").Value = "  -9.52%  "
$ws.Range("D16").Value = "'5.23"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "2.016.68"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "37.391.87"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'69.41"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "0.0₃0839"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "'5.13"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'227.41"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "'2.67"
$ws.Range("E23").Value = "  +6.55%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("D26").Value = "'163.21"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'8.89"
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("D28").Value = "'19.61"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -12.40%  "
$ws.Range("D30").Value = "'1.33"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").Value = "'0.0649"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'4.59"
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("D34").Value = "'4.48"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'5.18"
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("D41").Value = "'1.20"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").Value = "'0.0937"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("D43").Value = "'0.0212"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "1.396.59"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").Value = "'89.36"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "'15.60"
$ws.Range("E46").Value = "  -6.63%  "
$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").Value = "'7.05"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "2.206.21"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "'1.96"
$ws.Range("E51").Value = "  -2.58%  "
